$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" date column (C) for rows 2-6 from 45224 to 45233
$ws.Range("C2:C6").Value = 45233
